$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A2:AS2").ClearContents()
$ws.Range("B6").Select()
